$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3: "arch-gmail" -> "joong"
$ws.Range("D3").Value = "joong"

# Update fullname column (C2:C5): "Joong" -> "Joong A"
$ws.Range("C2").Value = "Joong A"
$ws.Range("C3").Value = "Joong A"
$ws.Range("C4").Value = "Joong A"
$ws.Range("C5").Value = "Joong A"

# Update phone numbers (E2, E3, E5): "081234562377" -> "08123456789"
# (leading apostrophe forces text storage so the leading zero is preserved)
$ws.Range("E2").Value = "'08123456789"
$ws.Range("E3").Value = "'08123456789"
$ws.Range("E5").Value = "'08123456789"

# Update D4: "archen1@gmail.com" -> "joong@gmail.com"
$ws.Range("D4").Value = "joong@gmail.com"

# Update selection to match new view state
$ws.Range("G16").Select() | Out-Null
